$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.178.00"
$ws.Range("E2").Value = "  -1.93%  "
$ws.Range("D3").Value = "2.177.37"
$ws.Range("E3").Value = "  -2.42%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.615"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "66.39"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -7.66%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.576"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "36.40"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -11.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0933"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.104"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.82"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.84%  "
$ws.Range("D15").Value = "2.503.42"
$ws.Range("E15").Value = "  -2.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.853"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.18%  "
$ws.Range("D18").Value = "2.183.23"
$ws.Range("E18").Value = "  -1.93%  "
$ws.Range("D19").Value = "41.101.85"
$ws.Range("E19").Value = "  -1.83%  "
$ws.Range("D20").Value = "0.0₃0948"
$ws.Range("E20").Value = "  -1.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.91%  "
$ws.Range("E22").Value = "  -2.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.25%  "
$ws.Range("E25").Value = "  -7.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  +5.62%  "
$ws.Range("E28").Value = "  -4.59%  "
$ws.Range("E29").Value = "  -3.77%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.00%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.02%  "
$ws.Range("E32").Value = "  -2.51%  "
$ws.Range("E33").Value = "  -1.47%  "
$ws.Range("E34").Value = "  +3.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0755"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.121"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.20%  "
$ws.Range("E37").Value = "  -4.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.99"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "24.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0306"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.40%  "
$ws.Range("E42").Value = "  +6.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.50"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "61.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.38%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.05%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.24%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.190"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.35%  "
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("E50").Value = "  -1.78%  "
$ws.Range("E51").Value = "  -3.98%  "
